$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the renamed station: "Spiegelgracht" -> "Rijksmuseum"
$ws.Range("A9").Value = "Rijksmuseum"

# Clear the lingering selection on D14, select A1 like a fresh view
$ws.Range("A1").Select()

# Autofit column A to the (now updated) content
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
